# Updated cryptos list on Sun Aug 20 18:19:51 UTC 2023 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) figures for each coin row,
# and reflects BabyDogeCoin/Aave swapping ranking positions (rows 45/46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and Volume(1h) (column E) updates
$ws.Range("D2").Value = "26.399.79"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.691.98"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "218.98"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "0.5488"
$ws.Range("E6").Value = "  +4.58%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "0.2728"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").Value = "0.06472"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").Value = "22.03"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("D11").Value = "0.07691"
$ws.Range("E11").Value = "  +3.11%  "
$ws.Range("D12").Value = "1.714.93"
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").Value = "4.549"
$ws.Range("D14").Value = "0.5846"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "0.000008407"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("D16").Value = "65.34"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").Value = "26.450.23"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "4.950"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "10.98"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "191.92"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "6.259"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "149.69"
$ws.Range("E24").Value = "  +3.38%  "
$ws.Range("D25").Value = "0.1329"
$ws.Range("E25").Value = "  +7.52%  "
$ws.Range("D26").Value = "7.898"
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("D27").Value = "15.75"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "0.06329"
$ws.Range("E28").Value = "  -4.68%  "
$ws.Range("D29").Value = "1.397"
$ws.Range("E29").Value = "  +2.96%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "3.601"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").Value = "3.602"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").Value = "1.685"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").Value = "1.045"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("D35").Value = "0.6165"
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("D36").Value = "2.409"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").Value = "6.225"
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("D39").Value = "1.119.77"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").Value = "0.01636"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("D41").Value = "0.8848"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "101.76"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("D44").Value = "1.842.03"
$ws.Range("E44").Value = "  +0.23%  "

# Row 45/46: BabyDogeCoin and Aave swapped positions (new ranking order)
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "57.53"
$ws.Range("E45").Value = "  +1.34%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.00000000107"
$ws.Range("E46").Value = "  -5.34%  "

$ws.Range("D47").Value = "8.227"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "0.05279"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").Value = "6.116"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").Value = "0.4303"
$ws.Range("E51").Value = "  -0.02%  "
